$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value could be misinterpreted as a number by Excel
# (e.g. "0.799", "1.00"). For these we temporarily force a Text number format
# so the value is stored as a string, then restore the default "Normal" style
# so the cell formatting matches the rest of the (unstyled) data cells.
function Set-TextValue($ws, $addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

$ws.Range('D2').Value = '46.664.14'
$ws.Range('E2').Value = '  +3.57%  '

$ws.Range('D3').Value = '2.257.27'
$ws.Range('E3').Value = '  -0.34%  '

$ws.Range('E4').Value = '  +0.09%  '

Set-TextValue $ws 'D5' '300.95'
$ws.Range('E5').Value = '  -0.59%  '

Set-TextValue $ws 'D6' '100.17'
$ws.Range('E6').Value = '  +5.15%  '

Set-TextValue $ws 'D7' '0.560'
$ws.Range('E7').Value = '  -0.80%  '

$ws.Range('E8').Value = '  +0.15%  '

Set-TextValue $ws 'D9' '0.511'
$ws.Range('E9').Value = '  +0.03%  '

Set-TextValue $ws 'D10' '35.51'
$ws.Range('E10').Value = '  +3.62%  '

Set-TextValue $ws 'D11' '0.0780'
$ws.Range('E11').Value = '  -1.36%  '

$ws.Range('E12').Value = '  -0.69%  '

$ws.Range('E13').Value = '  -0.80%  '

$ws.Range('D14').Value = '2.603.46'
$ws.Range('E14').Value = '  -0.26%  '

$ws.Range('D15').Value = '2.266.30'
$ws.Range('E15').Value = '  -0.22%  '

Set-TextValue $ws 'D16' '13.53'
$ws.Range('E16').Value = '  -0.68%  '

$ws.Range('B17').Value = 'Polygon'
$ws.Range('C17').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue $ws 'D17' '0.799'
$ws.Range('E17').Value = '  -0.09%  '

$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').Value = '46.640.95'
$ws.Range('E18').Value = '  +3.84%  '

Set-TextValue $ws 'D19' '12.92'
$ws.Range('E19').Value = '  -0.15%  '

$ws.Range('D20').Value = '0.0₃0926'
$ws.Range('E20').Value = '  +0.23%  '

Set-TextValue $ws 'D21' '5.88'
$ws.Range('E21').Value = '  -3.22%  '

Set-TextValue $ws 'D22' '65.25'
$ws.Range('E22').Value = '  -0.63%  '

Set-TextValue $ws 'D23' '248.90'
$ws.Range('E23').Value = '  +4.45%  '

$ws.Range('E24').Value = '  -2.15%  '

$ws.Range('E25').Value = '  +0.15%  '

$ws.Range('E26').Value = '  -0.97%  '

Set-TextValue $ws 'D27' '42.69'
$ws.Range('E27').Value = '  +2.78%  '

$ws.Range('E28').Value = '  -0.90%  '

Set-TextValue $ws 'D29' '9.68'
$ws.Range('E29').Value = '  +0.87%  '

$ws.Range('E30').Value = '  +1.34%  '

$ws.Range('E31').Value = '  +8.65%  '

Set-TextValue $ws 'D32' '146.53'
$ws.Range('E32').Value = '  -4.26%  '

Set-TextValue $ws 'D33' '5.42'
$ws.Range('E33').Value = '  -2.78%  '

$ws.Range('B34').Value = 'LidoDAOToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue $ws 'D34' '3.19'
$ws.Range('E34').Value = '  +7.81%  '

$ws.Range('B35').Value = 'Hedera'
$ws.Range('C35').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws 'D35' '0.0768'
$ws.Range('E35').Value = '  -2.68%  '

Set-TextValue $ws 'D36' '0.113'
$ws.Range('E36').Value = '  +8.47%  '

$ws.Range('E37').Value = '  -1.34%  '

Set-TextValue $ws 'D38' '16.18'
$ws.Range('E38').Value = '  +18.52%  '

Set-TextValue $ws 'D39' '1.71'
$ws.Range('E39').Value = '  -3.03%  '

Set-TextValue $ws 'D40' '3.84'
$ws.Range('E40').Value = '  -5.07%  '

Set-TextValue $ws 'D41' '0.0296'
$ws.Range('E41').Value = '  -6.09%  '

Set-TextValue $ws 'D42' '3.19'
$ws.Range('E42').Value = '  -2.10%  '

Set-TextValue $ws 'D43' '1.00'
$ws.Range('E43').Value = '  +0.06%  '

Set-TextValue $ws 'D44' '1.97'
$ws.Range('E44').Value = '  +0.92%  '

$ws.Range('D45').Value = '1.813.79'
$ws.Range('E45').Value = '  +3.96%  '

Set-TextValue $ws 'D46' '90.14'
$ws.Range('E46').Value = '  +18.54%  '

Set-TextValue $ws 'D47' '72.02'
$ws.Range('E47').Value = '  +1.35%  '

$ws.Range('E48').Value = '  -5.23%  '

Set-TextValue $ws 'D49' '4.80'
$ws.Range('E49').Value = '  +2.41%  '

Set-TextValue $ws 'D50' '93.58'
$ws.Range('E50').Value = '  -2.89%  '

$ws.Range('D51').Value = '2.480.39'
$ws.Range('E51').Value = '  -0.40%  '
